# Threat Alert Report update - 2026-01-18 01:00
# The oldest threat entry (row for 22-JAN-26) has rolled off the report;
# remove it and let the remaining rows shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (22-JAN-26 / SM-321 / Nile Air NP-119 ...).
# This shifts the 19-FEB-26 row up to row 2 and the 26-MAR-26 row up to row 3,
# carrying their existing values/styles with them.
$ws.Rows.Item(2).EntireRow.Delete()
